$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final roster data (player, position, team) for rows 2-19
$data = @(
    @("Josh Giddey",    "PG,SG,SF",    "Chicago Bulls"),
    @("Tyler Herro",    "PG,SG",       "Miami Heat"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Luka Doncic",    "PG,SG",       "Los Angeles Lakers"),
    @("Miles Bridges",  "SF,PF",       "Charlotte Hornets"),
    @("Brook Lopez",    "C",           "Milwaukee Bucks"),
    @("Zach Collins",   "PF,C",        "Chicago Bulls"),
    @("Christian Braun","SG,SF",       "Denver Nuggets"),
    @("Max Christie",   "SG,SF",       "Dallas Mavericks"),
    @("Isaiah Collier", "PG,SG",       "Utah Jazz"),
    @("Mikal Bridges",  "SG,SF,PF",    "New York Knicks"),
    @("De'Aaron Fox",   "PG,SG",       "San Antonio Spurs"),
    @("Ja Morant",      "PG",          "Memphis Grizzlies"),
    @("DeMar DeRozan",  "SF,PF",       "Sacramento Kings"),
    @("Evan Mobley",    "PF,C",        "Cleveland Cavaliers"),
    @("Nikola Vucevic", "PF,C",        "Chicago Bulls"),
    @("P.J. Washington","SF,PF",       "Dallas Mavericks"),
    @("Jaxson Hayes",   "PF,C",        "Los Angeles Lakers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
